# Taakverdeling.xlsx - add October/November task entries, restyle the
# previously-last row (13) so it no longer double-borders against the new
# row 14, and widen column A to fit the date column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New data rows 14-21
# ---------------------------------------------------------------------

# Row 14 - Ibrahim, 4 Oct 2024
$ws.Range("A14").Value = 45569
$ws.Range("B14").Value = "Ibrahim"
$ws.Range("C14").Value = 240
$ws.Range("D14").Value = "Account update & delete geimplementeerd"

# Row 15 - Ibrahim, 4 Oct 2024 (wraps to 2 lines)
$ws.Range("A15").Value = 45569
$ws.Range("B15").Value = "Ibrahim"
$ws.Range("C15").Value = 180
$ws.Range("D15").Value = "Account registratie pagina gemaakt de en functionaliteit hiervan gemaakt"

# Row 16 - Mark, 4 Nov 2024 (wraps to 2 lines)
$ws.Range("A16").Value = 45600
$ws.Range("B16").Value = "Mark"
$ws.Range("C16").Value = 360
$ws.Range("D16").Value = "Styling van het dashboard opnieuw bedacht en gemaakt"

# Row 17 - Mark, 4 Nov 2024 (wraps to 2 lines)
$ws.Range("A17").Value = 45600
$ws.Range("B17").Value = "Mark"
$ws.Range("C17").Value = 80
$ws.Range("D17").Value = "Verder gewerkt aan de navbar, want er waren meerdere problemen mee"

# Row 18 - Ibrahim, 4 Nov 2024 (wraps to 2 lines)
$ws.Range("A18").Value = 45600
$ws.Range("B18").Value = "Ibrahim"
$ws.Range("C18").Value = 80
$ws.Range("D18").Value = "Layout van de dashboard pagina opnieuw bedacht (ivm feedback stagebegeleider)"

# Row 19 - Ibrahim, 4 Nov 2024
$ws.Range("A19").Value = 45600
$ws.Range("B19").Value = "Ibrahim"
$ws.Range("C19").Value = 120
$ws.Range("D19").Value = "Login en registratie-pagina's afgemaakt"

# Row 20 - Ibrahim, 4 Nov 2024 (wraps to 2 lines)
$ws.Range("A20").Value = 45600
$ws.Range("B20").Value = "Ibrahim"
$ws.Range("C20").Value = 80
$ws.Range("D20").Value = "Database gegevens weergeven in de submenu van de navbar"

# Row 21 - Ibrahim, 4 Nov 2024 (wraps to 3 lines)
$ws.Range("A21").Value = 45600
$ws.Range("B21").Value = "Ibrahim"
$ws.Range("C21").Value = 80
$ws.Range("D21").Value = "Functionaliteiten van de account management aangepast zodat alleen admin accounts hier bij kunnen"

# ---------------------------------------------------------------------
# 2. Formatting: reuse the existing "Ibrahim" (row 2) and "Mark" (row 5)
#    banding so no redundant styles are minted.
# ---------------------------------------------------------------------

$ws.Range("A2:D2").Copy()
$ws.Range("A14:D15").PasteSpecial(-4122)   # xlPasteFormats - Ibrahim rows 14-15

$ws.Range("A5:D5").Copy()
$ws.Range("A16:D17").PasteSpecial(-4122)   # xlPasteFormats - Mark rows 16-17

$ws.Range("A2:D2").Copy()
$ws.Range("A18:D21").PasteSpecial(-4122)   # xlPasteFormats - Ibrahim rows 18-21

$excel.CutCopyMode = $false

# Row heights for the wrapped-text rows (2 or 3 visual lines).
$ws.Rows.Item(15).RowHeight = 28.8
$ws.Rows.Item(16).RowHeight = 28.8
$ws.Rows.Item(17).RowHeight = 28.8
$ws.Rows.Item(18).RowHeight = 28.8
$ws.Rows.Item(20).RowHeight = 28.8
$ws.Rows.Item(21).RowHeight = 43.2

# ---------------------------------------------------------------------
# 3. Row 13 is no longer the last row of the table: drop its bottom
#    border (keep left/right/top) so it doesn't double up against the
#    top border of the new row 14.
# ---------------------------------------------------------------------

$ws.Range("B13:D13").Borders.Item(9).LineStyle = 0   # xlEdgeBottom -> none

# ---------------------------------------------------------------------
# 4. Column A now needs to fit the widened date column.
# ---------------------------------------------------------------------

$ws.Columns.Item(1).ColumnWidth = 8.5

# ---------------------------------------------------------------------
# 5. Update the view: scroll so row 4 is at the top and select J16.
# ---------------------------------------------------------------------

$excel.ActiveWindow.ScrollRow = 4
$ws.Range("J16").Select()
